$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.194.05"
$ws.Range("E2").Value = "  +15.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.86"
$ws.Range("E3").Value = "  +8.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.07"
$ws.Range("E5").Value = "  +9.79%  "
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3731"
$ws.Range("E7").Value = "  +3.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3436"
$ws.Range("E8").Value = "  +8.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.07"
$ws.Range("E9").Value = "  +18.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.187"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07283"
$ws.Range("E11").Value = "  +7.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9972"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.59"
$ws.Range("E13").Value = "  +10.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.064"
$ws.Range("E14").Value = "  +7.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.768"
$ws.Range("E15").Value = "  +7.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.666.47"
$ws.Range("E16").Value = "  +8.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("E17").Value = "  +6.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9960"
$ws.Range("E18").Value = "  +3.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06716"
$ws.Range("E19").Value = "  +11.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.97"
$ws.Range("E20").Value = "  +14.57%  "
$ws.Range("E21").Value = "  +10.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.167"
$ws.Range("E22").Value = "  +9.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.03"
$ws.Range("E23").Value = "  +5.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.123.16"
$ws.Range("E24").Value = "  +14.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.401"
$ws.Range("E25").Value = "  +4.56%  "
$ws.Range("E26").Value = "  -8.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.669"
$ws.Range("E27").Value = "  +21.36%  "
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  +10.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.852.70"
$ws.Range("E30").Value = "  +9.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.38"
$ws.Range("E31").Value = "  +8.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.331"
$ws.Range("E32").Value = "  +22.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.018"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9863"
$ws.Range("E34").Value = "  +16.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.740"
$ws.Range("E35").Value = "  +16.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08424"
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.34"
$ws.Range("E37").Value = "  +15.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.953"
$ws.Range("E38").Value = "  +18.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06409"
$ws.Range("E39").Value = "  +9.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.337"
$ws.Range("E40").Value = "  +8.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.295"
$ws.Range("E41").Value = "  +6.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02347"
$ws.Range("E42").Value = "  +12.46%  "
$ws.Range("E43").Value = "  +11.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6131"
$ws.Range("E44").Value = "  +13.31%  "
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.806"
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.22"
$ws.Range("E47").Value = "  +6.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5979"
$ws.Range("E48").Value = "  +9.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.00"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.026"
$ws.Range("E50").Value = "  +8.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07119"
$ws.Range("E51").Value = "  +7.76%  "
